$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.88147134849534
$ws.Range("D2").Value = 7.049991988147958
$ws.Range("E2").Value = 10.13066848703873
$ws.Range("F2").Value = 36.69590696583991
$ws.Range("G2").Value = 3.693126488969042
$ws.Range("I2").Value = 29.09761494511422
$ws.Range("K2").Value = 14.56652409829984
$ws.Range("L2").Value = 10.60633411463625
$ws.Range("N2").Value = 21.50309568174834
$ws.Range("B3").Value = 13.78340319197619
$ws.Range("D3").Value = 7.069156225718231
$ws.Range("E3").Value = 10.0783698731568
$ws.Range("F3").Value = 36.51602977502817
$ws.Range("G3").Value = 3.696544905624772
$ws.Range("I3").Value = 29.14565157080244
$ws.Range("K3").Value = 14.27513527562455
$ws.Range("L3").Value = 10.56625742972427
$ws.Range("N3").Value = 21.55885782390413
$ws.Range("B4").Value = 13.72606961727194
$ws.Range("D4").Value = 7.081381407599278
$ws.Range("E4").Value = 10.04550668761286
$ws.Range("F4").Value = 36.41509646703552
$ws.Range("G4").Value = 3.698753338923227
$ws.Range("I4").Value = 29.17986127384013
$ws.Range("K4").Value = 14.09764170461958
$ws.Range("L4").Value = 10.54399577346527
$ws.Range("N4").Value = 21.59500616027011
$ws.Range("B5").Value = 13.70345229477011
$ws.Range("D5").Value = 7.086478975428452
$ws.Range("E5").Value = 10.03192625751189
$ws.Range("F5").Value = 36.37638463650018
$ws.Range("G5").Value = 3.699680931340195
$ws.Range("I5").Value = 29.19498558825541
$ws.Range("K5").Value = 14.02577858915098
$ws.Range("L5").Value = 10.53552011484274
$ws.Range("N5").Value = 21.61021773896412
$ws.Range("B6").Value = 13.69974242066616
$ws.Range("D6").Value = 7.087332424964996
$ws.Range("E6").Value = 10.02965987453842
$ws.Range("F6").Value = 36.37010341477026
$ws.Range("G6").Value = 3.699836629563858
$ws.Range("I6").Value = 29.1975683915852
$ws.Range("K6").Value = 14.01387720214392
$ws.Range("L6").Value = 10.53414891459684
$ws.Range("N6").Value = 21.61277266685266
$ws.Range("B7").Value = 13.7257615409248
$ws.Range("D7").Value = 7.081449686084505
$ws.Range("E7").Value = 10.04532429923298
$ws.Range("F7").Value = 36.41456455675853
$ws.Range("G7").Value = 3.698765736722168
$ws.Range("I7").Value = 29.18006045625562
$ws.Range("K7").Value = 14.09667049712827
$ws.Range("L7").Value = 10.54387904618952
$ws.Range("N7").Value = 21.59520936140406
$ws.Range("B8").Value = 13.84707326939161
$ws.Range("D8").Value = 7.056504991617736
$ws.Range("E8").Value = 10.11279227134436
$ws.Range("F8").Value = 36.63192868083959
$ws.Range("G8").Value = 3.694282490389518
$ws.Range("I8").Value = 29.11319818338339
$ws.Range("K8").Value = 14.46582395686886
$ws.Range("L8").Value = 10.59203286071293
$ws.Range("N8").Value = 21.52192632440241
$ws.Range("B9").Value = 14.1067885184385
$ws.Range("D9").Value = 7.011202491929843
$ws.Range("E9").Value = 10.2391106561565
$ws.Range("F9").Value = 37.13227551147099
$ws.Range("G9").Value = 3.686355208372451
$ws.Range("I9").Value = 29.0195803891065
$ws.Range("K9").Value = 15.19625617568116
$ws.Range("L9").Value = 10.70476310494054
$ws.Range("N9").Value = 21.39335084727692
$ws.Range("B10").Value = 14.3094478975894
$ws.Range("D10").Value = 6.980090145830084
$ws.Range("E10").Value = 10.32824305470194
$ws.Range("F10").Value = 37.54297195861698
$ws.Range("G10").Value = 3.68105156339189
$ws.Range("I10").Value = 28.97377286929963
$ws.Range("K10").Value = 15.73071139492243
$ws.Range("L10").Value = 10.79828179364392
$ws.Range("N10").Value = 21.30807751235159
$ws.Range("B11").Value = 14.40388719386475
$ws.Range("D11").Value = 6.966400992001969
$ws.Range("E11").Value = 10.36798895254996
$ws.Range("F11").Value = 37.73865462200325
$ws.Range("G11").Value = 3.678750459984327
$ws.Range("I11").Value = 28.95794308080005
$ws.Range("K11").Value = 15.97214063717177
$ws.Range("L11").Value = 10.8430366238581
$ws.Range("N11").Value = 21.2712733281525
$ws.Range("B12").Value = 14.43994374254728
$ws.Range("D12").Value = 6.961283479408377
$ws.Range("E12").Value = 10.3829239318424
$ws.Range("F12").Value = 37.81398037144022
$ws.Range("G12").Value = 3.677895028327223
$ws.Range("I12").Value = 28.95267026151445
$ws.Range("K12").Value = 16.06322332401079
$ws.Range("L12").Value = 10.86029193541111
$ws.Range("N12").Value = 21.2576218580381
$ws.Range("B13").Value = 14.43216568881116
$ws.Range("D13").Value = 6.962382686962801
$ws.Range("E13").Value = 10.3797125850107
$ws.Range("F13").Value = 37.79770393896588
$ws.Range("G13").Value = 3.678078553127479
$ws.Range("I13").Value = 28.9537737451327
$ws.Range("K13").Value = 16.04362370973094
$ws.Range("L13").Value = 10.8565621872168
$ws.Range("N13").Value = 21.26054925370796
$ws.Range("B14").Value = 14.40684785381706
$ws.Range("D14").Value = 6.965978645532483
$ws.Range("E14").Value = 10.36921999005755
$ws.Range("F14").Value = 37.74482745233804
$ws.Range("G14").Value = 3.678679764060528
$ws.Range("I14").Value = 28.95749481405166
$ws.Range("K14").Value = 15.97964144335802
$ws.Range("L14").Value = 10.84445012575505
$ws.Range("N14").Value = 21.27014449356198
$ws.Range("B15").Value = 14.39137742861559
$ws.Range("D15").Value = 6.968189891312902
$ws.Range("E15").Value = 10.36277785256206
$ws.Range("F15").Value = 37.71259712704882
$ws.Range("G15").Value = 3.679050096988941
$ws.Range("I15").Value = 28.95986808474918
$ws.Range("K15").Value = 15.94040321163545
$ws.Range("L15").Value = 10.83707087824555
$ws.Range("N15").Value = 21.27605902495693
$ws.Range("B16").Value = 14.30331885304101
$ws.Range("D16").Value = 6.980994062670447
$ws.Range("E16").Value = 10.32562938575362
$ws.Range("F16").Value = 37.53035774443493
$ws.Range("G16").Value = 3.681204182503597
$ws.Range("I16").Value = 28.97490828648125
$ws.Range("K16").Value = 15.71489083766873
$ws.Range("L16").Value = 10.79540059567006
$ws.Range("N16").Value = 21.31052271436391
$ws.Range("B17").Value = 14.24985325168945
$ws.Range("D17").Value = 6.988967529304331
$ws.Range("E17").Value = 10.30263433688146
$ws.Range("F17").Value = 37.42079415006615
$ws.Range("G17").Value = 3.682554147635565
$ws.Range("I17").Value = 28.9854187611886
$ws.Range("K17").Value = 15.57604072730103
$ws.Range("L17").Value = 10.77039661909961
$ws.Range("N17").Value = 21.33217381516234
$ws.Range("B18").Value = 14.21931465618379
$ws.Range("D18").Value = 6.993597357962758
$ws.Range("E18").Value = 10.28933285639111
$ws.Range("F18").Value = 37.35861259059087
$ws.Range("G18").Value = 3.683341117242228
$ws.Range("I18").Value = 28.99193544129608
$ws.Range("K18").Value = 15.496023568491
$ws.Range("L18").Value = 10.7562239319011
$ws.Range("N18").Value = 21.34481400372774
$ws.Range("B19").Value = 14.20901233818832
$ws.Range("D19").Value = 6.995172457932863
$ws.Range("E19").Value = 10.28481628785019
$ws.Range("F19").Value = 37.33770408648935
$ws.Range("G19").Value = 3.683609378827587
$ws.Range("I19").Value = 28.99422278319652
$ws.Range("K19").Value = 15.46890769263727
$ws.Range("L19").Value = 10.75146150162035
$ws.Range("N19").Value = 21.34912588849597
$ws.Range("B20").Value = 14.25552287773849
$ws.Range("D20").Value = 6.988114220540341
$ws.Range("E20").Value = 10.30509000155557
$ws.Range("F20").Value = 37.43237114712883
$ws.Range("G20").Value = 3.682409354967857
$ws.Range("I20").Value = 28.98425111145425
$ws.Range("K20").Value = 15.59083824444584
$ws.Range("L20").Value = 10.77303678550343
$ws.Range("N20").Value = 21.32984966186347
$ws.Range("B21").Value = 14.41427655701852
$ws.Range("D21").Value = 6.964920630290685
$ws.Range("E21").Value = 10.37230507192697
$ws.Range("F21").Value = 37.7603257105586
$ws.Range("G21").Value = 3.678502741811705
$ws.Range("I21").Value = 28.95638225241655
$ws.Range("K21").Value = 15.99844459796933
$ws.Range("L21").Value = 10.84799946963342
$ws.Range("N21").Value = 21.26731839327809
$ws.Range("B22").Value = 14.51973381503132
$ws.Range("D22").Value = 6.950148359220066
$ws.Range("E22").Value = 10.41555792632233
$ws.Range("F22").Value = 37.9817789097744
$ws.Range("G22").Value = 3.676042447708419
$ws.Range("I22").Value = 28.94237474905425
$ws.Range("K22").Value = 16.26280814843479
$ws.Range("L22").Value = 10.89877974929559
$ws.Range("N22").Value = 21.2281144602449
$ws.Range("B23").Value = 14.46330315046023
$ws.Range("D23").Value = 6.957997420419396
$ws.Range("E23").Value = 10.39253517738081
$ws.Range("F23").Value = 37.86295066778415
$ws.Range("G23").Value = 3.677347083459043
$ws.Range("I23").Value = 28.94946552451681
$ws.Range("K23").Value = 16.12192888983814
$ws.Range("L23").Value = 10.87151739471902
$ws.Range("N23").Value = 21.24888617415952
$ws.Range("B24").Value = 14.25295901622322
$ws.Range("D24").Value = 6.98849985879306
$ws.Range("E24").Value = 10.30398004885823
$ws.Range("F24").Value = 37.42713466619853
$ws.Range("G24").Value = 3.682474781908084
$ws.Range("I24").Value = 28.9847775293035
$ws.Range("K24").Value = 15.58414887246893
$ws.Range("L24").Value = 10.77184253492796
$ws.Range("N24").Value = 21.33089981133201
$ws.Range("B25").Value = 14.0343391140775
$ws.Range("D25").Value = 7.023074408779221
$ws.Range("E25").Value = 10.20558151874993
$ws.Range("F25").Value = 36.98918725560802
$ws.Range("G25").Value = 3.688407877182786
$ws.Range("I25").Value = 29.04087955769291
$ws.Range("K25").Value = 14.99860700921628
$ws.Range("L25").Value = 10.67235227550335
$ws.Range("N25").Value = 21.42651750364235

Write-Output "Updated 216 cells (rows 2-25, columns B,D,E,F,G,I,K,L,N)"
